$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5850
$ws.Range("I62").Value = 5800
$ws.Range("K62").Value = 5800
$ws.Range("M62").Value = -5176

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2459
$ws.Range("I64").Value = 2431.6667
$ws.Range("J64").Value = 2500
$ws.Range("K64").Value = 2431.6667
$ws.Range("L64").Value = 2500
$ws.Range("M64").Value = -2183.6667
$ws.Range("N64").Value = -2996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5850
$ws.Range("I65").Value = 5800
$ws.Range("K65").Value = 29000
$ws.Range("M65").Value = -25880

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 2459
$ws.Range("I67").Value = 2431.6667
$ws.Range("J67").Value = 2500
$ws.Range("K67").Value = 2431.6667
$ws.Range("L67").Value = 2500
$ws.Range("M67").Value = -1573.6667
$ws.Range("N67").Value = -4216

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 907.4
$ws.Range("I80").Value = 733.3333
$ws.Range("K80").Value = 2199.9999
$ws.Range("M80").Value = -1201.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 907.4
$ws.Range("I83").Value = 733.3333
$ws.Range("K83").Value = 6599.9997
$ws.Range("M83").Value = -1607.9997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5166.6665
$ws.Range("I100").Value = 5166.6665
$ws.Range("K100").Value = 5166.6665
$ws.Range("M100").Value = -4625.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2969.4443
$ws.Range("I137").Value = 1909.3334
$ws.Range("J137").Value = 3499.5
$ws.Range("K137").Value = 5728.0002
$ws.Range("L137").Value = 10498.5
$ws.Range("M137").Value = -3178.0002
$ws.Range("N137").Value = -15598.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 65.5
$ws.Range("I5").Value = 67.333336
$ws.Range("K5").Value = 67.333336
$ws.Range("M5").Value = 44.666664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2366.6758
$ws.Range("I32").Value = 2321.3057
$ws.Range("K32").Value = 2321.3057
$ws.Range("M32").Value = -2034.3057

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3700
$ws.Range("I63").Value = 3700
$ws.Range("K63").Value = 3700
$ws.Range("M63").Value = -3014

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3700
$ws.Range("I66").Value = 3700
$ws.Range("K66").Value = 18500
$ws.Range("M66").Value = -15068

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 65.5
$ws.Range("I4").Value = 67.333336
$ws.Range("K4").Value = 67.333336
$ws.Range("M4").Value = 47.666664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 655.8
$ws.Range("I22").Value = 740
$ws.Range("J22").Value = 599.6667
$ws.Range("K22").Value = 740
$ws.Range("L22").Value = 599.6667
$ws.Range("M22").Value = -567
$ws.Range("N22").Value = -945.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 52500
$ws.Range("J63").Value = 52500
$ws.Range("L63").Value = 52500
$ws.Range("N63").Value = -53872

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H66").Value = 52500
$ws.Range("J66").Value = 52500
$ws.Range("L66").Value = 157500
$ws.Range("N66").Value = -164364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2666.3333
$ws.Range("I105").Value = 2666.3333
$ws.Range("K105").Value = 2666.3333
$ws.Range("M105").Value = -919.3332999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5320
$ws.Range("I31").Value = 1793.8462
$ws.Range("J31").Value = 9140
$ws.Range("K31").Value = 1793.8462
$ws.Range("L31").Value = 9140
$ws.Range("M31").Value = -1498.8462
$ws.Range("N31").Value = -9730

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5320
$ws.Range("I34").Value = 1793.8462
$ws.Range("J34").Value = 9140
$ws.Range("K34").Value = 1793.8462
$ws.Range("L34").Value = 9140
$ws.Range("M34").Value = -1591.8462
$ws.Range("N34").Value = -9544

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 30046
$ws.Range("J97").Value = 30046
$ws.Range("L97").Value = 30046
$ws.Range("N97").Value = -32028

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 95605680
$ws.Range("I4").Value = 17499798
$ws.Range("K4").Value = 52499394
$ws.Range("M4").Value = -52499282

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2000
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 5524.25
$ws.Range("J107").Value = 1100
$ws.Range("L107").Value = 3300
$ws.Range("N107").Value = -7140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 150
$ws.Range("I2").Value = 200
$ws.Range("K2").Value = 200
$ws.Range("M2").Value = -87

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2248.5
$ws.Range("I126").Value = 2499
$ws.Range("J126").Value = 1998
$ws.Range("K126").Value = 7497
$ws.Range("L126").Value = 5994
$ws.Range("M126").Value = -5027
$ws.Range("N126").Value = -10934

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2703
$ws.Range("I46").Value = 2703
$ws.Range("K46").Value = 2703
$ws.Range("M46").Value = -2515

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3606.9285
$ws.Range("I132").Value = 2278.6667
$ws.Range("J132").Value = 5997.8
$ws.Range("K132").Value = 6836.000100000001
$ws.Range("L132").Value = 17993.4
$ws.Range("M132").Value = -4306.000100000001
$ws.Range("N132").Value = -23053.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 33165
$ws.Range("J54").Value = 39198
$ws.Range("L54").Value = 39198
$ws.Range("N54").Value = -40238

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6074.3335
$ws.Range("I81").Value = 6074.3335
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 12148.667
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -11087.667
$ws.Range("N81").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 6074.3335
$ws.Range("I84").Value = 6074.3335
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 60743.335
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -55439.335
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2590.9033
$ws.Range("I132").Value = 2410.6365
$ws.Range("J132").Value = 3031.5557
$ws.Range("K132").Value = 7231.9095
$ws.Range("L132").Value = 9094.667099999999
$ws.Range("M132").Value = -4701.9095
$ws.Range("N132").Value = -14154.6671
